$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 30000
$ws.Range("J10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("N10").Value = -30586

$ws.Range("H112").Value = 27778798
$ws.Range("I112").Value = 725
$ws.Range("J112").Value = 35715390
$ws.Range("K112").Value = 2175
$ws.Range("L112").Value = 107146170
$ws.Range("M112").Value = -1067
$ws.Range("N112").Value = -107148386

$ws.Range("H113").Value = 12503130
$ws.Range("I113").Value = 3538.125
$ws.Range("J113").Value = 62501500
$ws.Range("K113").Value = 3538.125
$ws.Range("L113").Value = 62501500
$ws.Range("M113").Value = -284.125
$ws.Range("N113").Value = -62508008

$ws.Range("H116").Value = 7325.8
$ws.Range("I116").Value = 11681
$ws.Range("J116").Value = 2970.6
$ws.Range("K116").Value = 11681
$ws.Range("L116").Value = 2970.6
$ws.Range("M116").Value = -8239
$ws.Range("N116").Value = -9854.6

$ws.Range("H125").Value = 4377
$ws.Range("I125").Value = 9177.333000000001
$ws.Range("J125").Value = 3348.3572
$ws.Range("K125").Value = 82595.997
$ws.Range("L125").Value = 30135.2148
$ws.Range("M125").Value = -80135.997
$ws.Range("N125").Value = -35055.2148

$ws.Range("H129").Value = 906.41895
$ws.Range("I129").Value = 578.26666
$ws.Range("K129").Value = 1734.79998
$ws.Range("M129").Value = 3265.20002

$ws.Range("H132").Value = 10102565
$ws.Range("I132").Value = 1262.4166
$ws.Range("K132").Value = 3787.2498
$ws.Range("M132").Value = -1257.2498

$ws.Range("H135").Value = 1488.7142
$ws.Range("I135").Value = 1185.9762
$ws.Range("K135").Value = 10673.7858
$ws.Range("M135").Value = -8138.785800000001

$ws.Range("H137").Value = 1519.4524
$ws.Range("I137").Value = 1266.1111
$ws.Range("J137").Value = 1975.4667
$ws.Range("K137").Value = 3798.3333
$ws.Range("L137").Value = 5926.4001
$ws.Range("M137").Value = -1248.3333
$ws.Range("N137").Value = -11026.4001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10105.917
$ws.Range("I45").Value = 12107.889
$ws.Range("J45").Value = 4100
$ws.Range("K45").Value = 12107.889
$ws.Range("L45").Value = 4100
$ws.Range("M45").Value = -11730.889
$ws.Range("N45").Value = -4854

$ws.Range("H61").Value = 4669.484
$ws.Range("I61").Value = 4871.1035
$ws.Range("K61").Value = 4871.1035
$ws.Range("M61").Value = -4659.1035

$ws.Range("H97").Value = 1164.3334
$ws.Range("I97").Value = 1014.3125
$ws.Range("J97").Value = 1644.4
$ws.Range("K97").Value = 1014.3125
$ws.Range("L97").Value = 1644.4
$ws.Range("M97").Value = -518.3125
$ws.Range("N97").Value = -2636.4

$ws.Range("H132").Value = 2439.42
$ws.Range("I132").Value = 1825.5143
$ws.Range("J132").Value = 3871.8667
$ws.Range("K132").Value = 5476.5429
$ws.Range("L132").Value = 11615.6001
$ws.Range("M132").Value = -2946.5429
$ws.Range("N132").Value = -16675.6001

$ws.Range("H136").Value = 4669.484
$ws.Range("I136").Value = 4871.1035
$ws.Range("K136").Value = 14613.3105
$ws.Range("M136").Value = -12063.3105

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 333334660
$ws.Range("I99").Value = 333334660
$ws.Range("K99").Value = 333334660
$ws.Range("M99").Value = -333333162

$ws.Range("H115").Value = 59684
$ws.Range("J115").Value = 59684
$ws.Range("L115").Value = 59684
$ws.Range("N115").Value = -62818

$ws.Range("H134").Value = 3011.791
$ws.Range("I134").Value = 3397.2917
$ws.Range("J134").Value = 2037.8948
$ws.Range("K134").Value = 10191.8751
$ws.Range("L134").Value = 6113.6844
$ws.Range("M134").Value = -7656.875100000001
$ws.Range("N134").Value = -11183.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5926414
$ws.Range("I31").Value = 1575.7551
$ws.Range("J31").Value = 13566337
$ws.Range("K31").Value = 1575.7551
$ws.Range("L31").Value = 13566337
$ws.Range("M31").Value = -1280.7551
$ws.Range("N31").Value = -13566927

$ws.Range("H34").Value = 5926414
$ws.Range("I34").Value = 1575.7551
$ws.Range("J34").Value = 13566337
$ws.Range("K34").Value = 1575.7551
$ws.Range("L34").Value = 13566337
$ws.Range("M34").Value = -1373.7551
$ws.Range("N34").Value = -13566741

$ws.Range("I58").Value = 3704262
$ws.Range("J58").Value = 2213.0557
$ws.Range("K58").Value = 3704262
$ws.Range("L58").Value = 2213.0557
$ws.Range("M58").Value = -3704059
$ws.Range("N58").Value = -2619.0557

$ws.Range("H62").Value = 5178.9287
$ws.Range("I62").Value = 5192.6924
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 5192.6924
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -4568.6924
$ws.Range("N62").Value = -6248

$ws.Range("H63").Value = 42300
$ws.Range("J63").Value = 42300
$ws.Range("L63").Value = 42300
$ws.Range("N63").Value = -43672

$ws.Range("H65").Value = 5178.9287
$ws.Range("I65").Value = 5192.6924
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 25963.462
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -22843.462
$ws.Range("N65").Value = -31240

$ws.Range("H66").Value = 42300
$ws.Range("J66").Value = 42300
$ws.Range("L66").Value = 126900
$ws.Range("N66").Value = -133764

$ws.Range("H68").Value = 25283.166
$ws.Range("J68").Value = 28339.8
$ws.Range("L68").Value = 28339.8
$ws.Range("N68").Value = -29837.8

$ws.Range("H71").Value = 25283.166
$ws.Range("J71").Value = 28339.8
$ws.Range("L71").Value = 85019.39999999999
$ws.Range("N71").Value = -92507.39999999999

$ws.Range("H99").Value = 6953884
$ws.Range("I99").Value = 10182.909
$ws.Range("J99").Value = 17865414
$ws.Range("K99").Value = 10182.909
$ws.Range("L99").Value = 17865414
$ws.Range("M99").Value = -8684.909
$ws.Range("N99").Value = -17868410

$ws.Range("H126").Value = 6953884
$ws.Range("I126").Value = 10182.909
$ws.Range("J126").Value = 17865414
$ws.Range("K126").Value = 30548.727
$ws.Range("L126").Value = 53596242
$ws.Range("M126").Value = -28078.727
$ws.Range("N126").Value = -53601182

$ws.Range("H132").Value = 3573084.8
$ws.Range("I132").Value = 4546883.5
$ws.Range("J132").Value = 2489.75
$ws.Range("K132").Value = 13640650.5
$ws.Range("L132").Value = 7469.25
$ws.Range("M132").Value = -13638120.5
$ws.Range("N132").Value = -12529.25

$ws.Range("H134").Value = 7938777.5
$ws.Range("I134").Value = 13891883
$ws.Range("K134").Value = 41675649
$ws.Range("M134").Value = -41673114

$ws.Range("I136").Value = 3704262
$ws.Range("J136").Value = 2213.0557
$ws.Range("K136").Value = 11112786
$ws.Range("L136").Value = 6639.1671
$ws.Range("M136").Value = -11110236
$ws.Range("N136").Value = -11739.1671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 9251
$ws.Range("J93").Value = 9251
$ws.Range("L93").Value = 9251
$ws.Range("N93").Value = -12995

$ws.Range("H113").Value = 62501444
$ws.Range("I113").Value = 125001130
$ws.Range("J113").Value = 1764.125
$ws.Range("K113").Value = 125001130
$ws.Range("L113").Value = 1764.125
$ws.Range("M113").Value = -124998960
$ws.Range("N113").Value = -6104.125

$ws.Range("H132").Value = 4275147
$ws.Range("I132").Value = 4506171.5
$ws.Range("J132").Value = 1195
$ws.Range("K132").Value = 13518514.5
$ws.Range("L132").Value = 3585
$ws.Range("M132").Value = -13515984.5
$ws.Range("N132").Value = -8645

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 10888021
$ws.Range("I132").Value = 13012064
$ws.Range("K132").Value = 39036192
$ws.Range("M132").Value = -39033662

$ws.Range("H136").Value = 5689.912
$ws.Range("I136").Value = 4081.2954
$ws.Range("J136").Value = 11134.462
$ws.Range("K136").Value = 12243.8862
$ws.Range("L136").Value = 33403.386
$ws.Range("M136").Value = -9693.886200000001
$ws.Range("N136").Value = -38503.386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 832.38184
$ws.Range("I132").Value = 601.15
$ws.Range("J132").Value = 1449
$ws.Range("K132").Value = 1803.45
$ws.Range("L132").Value = 4347
$ws.Range("M132").Value = 726.5500000000002
$ws.Range("N132").Value = -9407

$ws.Range("H136").Value = 13074542
$ws.Range("I136").Value = 3410.7666
$ws.Range("J136").Value = 31747586
$ws.Range("K136").Value = 10232.2998
$ws.Range("L136").Value = 95242758
$ws.Range("M136").Value = -7682.299800000001
$ws.Range("N136").Value = -95247858
